# Add a second "Logging" sheet ahead of the existing timesheet, which is
# renamed to "Timesheet". Mirrors commit "update templates with second
# 'Logging' timesheet".

$wb = $excel.ActiveWorkbook

# The single pre-existing sheet becomes "Timesheet".
$timesheet = $wb.Worksheets.Item(1)
$timesheet.Name = "Timesheet"

# New "Logging" sheet, inserted before "Timesheet" so it becomes the first
# (leftmost) tab.
$logging = $wb.Worksheets.Add($timesheet)
$logging.Name = "Logging"

# Small carry-over bookkeeping table used by the logging sheet.
$logging.Range("B1").Value = "carryover"
$logging.Range("A2").Value = "row"
$logging.Range("B2").Value = 35
$logging.Range("A3").Value = "column"
$logging.Range("B3").Value = 10

# Re-fetch the Timesheet worksheet by name: after inserting/reordering
# sheets, stale worksheet references don't reliably propagate selection /
# activation, so grab a fresh handle before selecting + activating it.
$timesheet = $wb.Worksheets.Item("Timesheet")

# Move the selection on the timesheet tab.
$timesheet.Range("J35").Select() | Out-Null

# Timesheet is the active/visible tab when the workbook is opened.
$timesheet.Activate() | Out-Null
